$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-slot labels in column C (rows 2 and 3) with new values.
$ws.Range("C2").Value = "2:55-3:0"
$ws.Range("C3").Value = "3:0-3:5"

# Move the active selection to B11 (single cell), matching the final cursor
# position left in the worksheet after the edit.
[void]$ws.Range("B11").Select()
